# Update NATMI LR-pair TPM output (Adam9-Itgb5) with newly recomputed TPM values.
# Only the "ECs" sending/target cluster's underlying ligand/receptor expression
# values actually changed; every other column in the sheet (specificity scores,
# edge weights, edge specificity scores) is derived from those base values and
# is therefore recomputed here so the whole table stays internally consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 10

# --- 1. Write the new base (directly-updated) values -----------------------
# Ligand average / total expression value for the "ECs" sending cluster
# (column G / H), and Receptor average / total expression value for the
# "ECs" target cluster (column M / N).
$newLigandAvg  = 5.685592333333333   # ECs ligand average expression value (G)
$newLigandTot  = 17.056777           # ECs ligand total expression value   (H)
$newReceptAvg  = 9.101794333333332   # ECs receptor average expression value (M)
$newReceptTot  = 27.305383           # ECs receptor total expression value   (N)

for ($r = 2; $r -le $lastRow; $r++) {
    $sendCluster = $ws.Cells.Item($r, 1).Value()
    $targCluster = $ws.Cells.Item($r, 4).Value()

    if ($sendCluster -eq "ECs") {
        $ws.Cells.Item($r, 7).Value = $newLigandAvg
        $ws.Cells.Item($r, 8).Value = $newLigandTot
    }
    if ($targCluster -eq "ECs") {
        $ws.Cells.Item($r, 13).Value = $newReceptAvg
        $ws.Cells.Item($r, 14).Value = $newReceptTot
    }
}

# --- 2. Recompute the derived specificity / edge-weight columns ------------
# Ligand derived specificity = ligand value for this row's sending cluster
# divided by the sum of the ligand values across all sending clusters
# (columns I and J use average / total respectively). Likewise for the
# receptor-side derived specificity (columns O and P).
$sumLigandAvg = 0.0
$sumLigandTot = 0.0
$sumReceptAvg = 0.0
$sumReceptTot = 0.0

$seenSend = @{}
$seenTarg = @{}

for ($r = 2; $r -le $lastRow; $r++) {
    $sendCluster = $ws.Cells.Item($r, 1).Value()
    $targCluster = $ws.Cells.Item($r, 4).Value()

    if (-not $seenSend.ContainsKey($sendCluster)) {
        $seenSend[$sendCluster] = $true
        $sumLigandAvg += $ws.Cells.Item($r, 7).Value()
        $sumLigandTot += $ws.Cells.Item($r, 8).Value()
    }
    if (-not $seenTarg.ContainsKey($targCluster)) {
        $seenTarg[$targCluster] = $true
        $sumReceptAvg += $ws.Cells.Item($r, 13).Value()
        $sumReceptTot += $ws.Cells.Item($r, 14).Value()
    }
}

$sumEdgeAvg = 0.0
$sumEdgeTot = 0.0

for ($r = 2; $r -le $lastRow; $r++) {
    $sendCluster = $ws.Cells.Item($r, 1).Value()
    $targCluster = $ws.Cells.Item($r, 4).Value()

    $ligAvg = $ws.Cells.Item($r, 7).Value()
    $ligTot = $ws.Cells.Item($r, 8).Value()
    $recAvg = $ws.Cells.Item($r, 13).Value()
    $recTot = $ws.Cells.Item($r, 14).Value()

    $ws.Cells.Item($r, 9).Value  = $ligAvg / $sumLigandAvg   # I: ligand specificity (avg)
    $ws.Cells.Item($r, 10).Value = $ligTot / $sumLigandTot   # J: ligand specificity (total)
    $ws.Cells.Item($r, 15).Value = $recAvg / $sumReceptAvg   # O: receptor specificity (avg)
    $ws.Cells.Item($r, 16).Value = $recTot / $sumReceptTot   # P: receptor specificity (total)

    # Only re-derive the edge weight (Q/R = ligand * receptor value) for rows whose
    # underlying ligand or receptor value actually changed (i.e. involves "ECs"),
    # leaving untouched rows bit-for-bit identical to their original values.
    if ($sendCluster -eq "ECs" -or $targCluster -eq "ECs") {
        $edgeAvg = $ligAvg * $recAvg                          # Q: edge average expression weight
        $edgeTot = $ligTot * $recTot                          # R: edge total expression weight
        $ws.Cells.Item($r, 17).Value = $edgeAvg
        $ws.Cells.Item($r, 18).Value = $edgeTot
    }

    $sumEdgeAvg += $ws.Cells.Item($r, 17).Value()
    $sumEdgeTot += $ws.Cells.Item($r, 18).Value()
}

for ($r = 2; $r -le $lastRow; $r++) {
    $edgeAvg = $ws.Cells.Item($r, 17).Value()
    $edgeTot = $ws.Cells.Item($r, 18).Value()

    $ws.Cells.Item($r, 19).Value = $edgeAvg / $sumEdgeAvg   # S: edge average expression derived specificity
    $ws.Cells.Item($r, 20).Value = $edgeTot / $sumEdgeTot   # T: edge total expression derived specificity
}
